# This workbook contains one sheet with a weekly price table for
# "Poroto verde" at "Terminal La Palmera de La Serena".
# The edit inserts one new data row (a new weekly observation) right
# before the existing row 398, which pushes every following row down
# by one (old row 398 becomes 399, ..., old row 490 becomes 491).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 398; this shifts rows 398..490 down to 399..491
# and automatically grows the sheet dimension from A1:R490 to A1:R491.
$ws.Rows.Item(398).Insert()

# Populate the newly inserted row 398 with the new observation.
$ws.Range("A398").Value = 8
$ws.Range("B398").Value = "Terminal La Palmera de La Serena"
$ws.Range("C398").Value = "Coquimbo"
$ws.Range("D398").Value = 45275
$ws.Range("E398").Value = 4
$ws.Range("F398").Value = 100112031
$ws.Range("G398").Value = "Poroto verde"
$ws.Range("H398").Value = "Sin especificar"
$ws.Range("I398").Value = "Primera"
$ws.Range("J398").Value = 400
$ws.Range("K398").Value = 29000
$ws.Range("L398").Value = 30000
$ws.Range("M398").Value = 29500
$ws.Range("N398").Value = "$/malla 25 kilos"
$ws.Range("O398").Value = "Región de Arica y Parinacota"
$ws.Range("P398").Value = 1180
$ws.Range("Q398").Value = 25
$ws.Range("R398").Value = "Hortaliza"
